# Insert a new data row at row 3 (pushes existing rows 3..75 down to 4..76)
# and populate it with a new price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(3, 3).Value = "Maule"
$ws.Cells.Item(3, 4).Value = 44515
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100103
$ws.Cells.Item(3, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(3, 9).Value = 100103001
$ws.Cells.Item(3, 10).Value = "Cereza"
$ws.Cells.Item(3, 11).Value = "Royal Dawn"
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 60
$ws.Cells.Item(3, 14).Value = 20000
$ws.Cells.Item(3, 15).Value = 20000
$ws.Cells.Item(3, 16).Value = 20000
$ws.Cells.Item(3, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(3, 19).Value = 2000
$ws.Cells.Item(3, 20).Value = 10
